$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.296496629714966
$ws.Range("B1").Value = 1.116909027099609
$ws.Range("C1").Value = 3.172342538833618
$ws.Range("D1").Value = 3.100940227508545
$ws.Range("E1").Value = 0.9243699908256531
